$wb = $excel.ActiveWorkbook

# Overview sheet: G2 timestamp update
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-20 13:05:44"

# zh-cn sheet: H2 and K2 timestamp updates
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-20 13:05:40"
$wsZhCn.Range("K2").Value = "2016-08-20 13:05:57"

# de-de sheet: K2 timestamp update
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-20 13:06:08"
